$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2296.7334
$ws.Range("I33").Value = 314.27274
$ws.Range("K33").Value = 314.27274
$ws.Range("M33").Value = -85.27274
$ws.Range("H74").Value = 6446.75
$ws.Range("I74").Value = 6153.4287
$ws.Range("K74").Value = 6153.4287
$ws.Range("M74").Value = -5217.4287
$ws.Range("H77").Value = 6446.75
$ws.Range("I77").Value = 6153.4287
$ws.Range("K77").Value = 30767.1435
$ws.Range("M77").Value = -26087.1435
$ws.Range("H88").Value = 5635.4
$ws.Range("J88").Value = 3089.5
$ws.Range("L88").Value = 3089.5
$ws.Range("N88").Value = -3901.5
$ws.Range("H91").Value = 5635.4
$ws.Range("J91").Value = 3089.5
$ws.Range("L91").Value = 3089.5
$ws.Range("N91").Value = -5897.5
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("M94").Value = -549
$ws.Range("H98").Value = 857.1739
$ws.Range("I98").Value = 607
$ws.Range("K98").Value = 607
$ws.Range("M98").Value = 891
$ws.Range("H100").Value = 1839
$ws.Range("I100").Value = 772.6
$ws.Range("K100").Value = 772.6
$ws.Range("M100").Value = -231.6
$ws.Range("H122").Value = 857.1739
$ws.Range("I122").Value = 607
$ws.Range("K122").Value = 1821
$ws.Range("M122").Value = 629
$ws.Range("H132").Value = 2548.5908
$ws.Range("I132").Value = 2584.238
$ws.Range("K132").Value = 7752.714
$ws.Range("M132").Value = -5222.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5670.0713
$ws.Range("I32").Value = 2506.5938
$ws.Range("K32").Value = 2506.5938
$ws.Range("M32").Value = -2219.5938
$ws.Range("H88").Value = 1668.8889
$ws.Range("I88").Value = 1140
$ws.Range("J88").Value = 2092
$ws.Range("K88").Value = 1140
$ws.Range("L88").Value = 2092
$ws.Range("M88").Value = -734
$ws.Range("N88").Value = -2904
$ws.Range("H91").Value = 1668.8889
$ws.Range("I91").Value = 1140
$ws.Range("J91").Value = 2092
$ws.Range("K91").Value = 1140
$ws.Range("L91").Value = 2092
$ws.Range("M91").Value = 264
$ws.Range("N91").Value = -4900
$ws.Range("H97").Value = 1561.6875
$ws.Range("I97").Value = 1561.6875
$ws.Range("K97").Value = 1561.6875
$ws.Range("M97").Value = -1065.6875
$ws.Range("H131").Value = 83777
$ws.Range("J131").Value = 83777
$ws.Range("L131").Value = 83777
$ws.Range("N131").Value = -93857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2003.1538
$ws.Range("I86").Value = 1360.5555
$ws.Range("J86").Value = 3449
$ws.Range("K86").Value = 1360.5555
$ws.Range("L86").Value = 3449
$ws.Range("M86").Value = -237.5554999999999
$ws.Range("N86").Value = -5695
$ws.Range("H89").Value = 2003.1538
$ws.Range("I89").Value = 1360.5555
$ws.Range("J89").Value = 3449
$ws.Range("K89").Value = 6802.7775
$ws.Range("L89").Value = 17245
$ws.Range("M89").Value = -1186.7775
$ws.Range("N89").Value = -28477
$ws.Range("H99").Value = 2437.3809
$ws.Range("I99").Value = 2158.1177
$ws.Range("K99").Value = 2158.1177
$ws.Range("M99").Value = -660.1176999999998
$ws.Range("H105").Value = 2603.7368
$ws.Range("I105").Value = 2553
$ws.Range("K105").Value = 2553
$ws.Range("M105").Value = -806

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 202.14285
$ws.Range("I7").Value = 178
$ws.Range("J7").Value = 234.33333
$ws.Range("K7").Value = 178
$ws.Range("L7").Value = 234.33333
$ws.Range("M7").Value = -65
$ws.Range("N7").Value = -460.33333
$ws.Range("H58").Value = 3364.2856
$ws.Range("I58").Value = 1137.5
$ws.Range("J58").Value = 6333.3335
$ws.Range("K58").Value = 1137.5
$ws.Range("L58").Value = 6333.3335
$ws.Range("M58").Value = -934.5
$ws.Range("N58").Value = -6739.3335
$ws.Range("H60").Value = 32041.5
$ws.Range("I60").Value = 21750
$ws.Range("J60").Value = 42333
$ws.Range("K60").Value = 21750
$ws.Range("L60").Value = 42333
$ws.Range("M60").Value = -21239
$ws.Range("N60").Value = -43355
$ws.Range("H86").Value = 10007
$ws.Range("I86").Value = 10007
$ws.Range("K86").Value = 10007
$ws.Range("M86").Value = -8884
$ws.Range("H89").Value = 10007
$ws.Range("I89").Value = 10007
$ws.Range("K89").Value = 50035
$ws.Range("M89").Value = -44419
$ws.Range("H93").Value = 167984.5
$ws.Range("I93").Value = 10407
$ws.Range("J93").Value = 199500
$ws.Range("K93").Value = 10407
$ws.Range("L93").Value = 199500
$ws.Range("M93").Value = -8535
$ws.Range("N93").Value = -203244
$ws.Range("H105").Value = 2181
$ws.Range("I105").Value = 2252
$ws.Range("K105").Value = 2252
$ws.Range("M105").Value = -505
$ws.Range("H107").Value = 2179.1667
$ws.Range("I107").Value = 1940
$ws.Range("J107").Value = 3375
$ws.Range("K107").Value = 1940
$ws.Range("L107").Value = 3375
$ws.Range("M107").Value = -20
$ws.Range("N107").Value = -7215
$ws.Range("H132").Value = 2597.2273
$ws.Range("I132").Value = 2050
$ws.Range("J132").Value = 6865.6
$ws.Range("K132").Value = 6150
$ws.Range("L132").Value = 20596.8
$ws.Range("M132").Value = -3620
$ws.Range("N132").Value = -25656.8
$ws.Range("H136").Value = 3364.2856
$ws.Range("I136").Value = 1137.5
$ws.Range("J136").Value = 6333.3335
$ws.Range("K136").Value = 3412.5
$ws.Range("L136").Value = 19000.0005
$ws.Range("M136").Value = -862.5
$ws.Range("N136").Value = -24100.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10008.533
$ws.Range("I80").Value = 8482.842000000001
$ws.Range("J80").Value = 12643.818
$ws.Range("K80").Value = 8482.842000000001
$ws.Range("L80").Value = 12643.818
$ws.Range("M80").Value = -7484.842000000001
$ws.Range("N80").Value = -14639.818
$ws.Range("H83").Value = 10008.533
$ws.Range("I83").Value = 8482.842000000001
$ws.Range("J83").Value = 12643.818
$ws.Range("K83").Value = 42414.21000000001
$ws.Range("L83").Value = 63219.09
$ws.Range("M83").Value = -37422.21000000001
$ws.Range("N83").Value = -73203.09
$ws.Range("H97").Value = 884.9545000000001
$ws.Range("I97").Value = 399.08334
$ws.Range("J97").Value = 1067.1562
$ws.Range("K97").Value = 399.08334
$ws.Range("L97").Value = 1067.1562
$ws.Range("M97").Value = 96.91665999999998
$ws.Range("N97").Value = -2059.1562
$ws.Range("H102").Value = 2750.1187
$ws.Range("I102").Value = 2315.9148
$ws.Range("K102").Value = 2315.9148
$ws.Range("M102").Value = -693.9148
$ws.Range("H113").Value = 3161.3928
$ws.Range("I113").Value = 2469.9
$ws.Range("J113").Value = 4890.125
$ws.Range("K113").Value = 2469.9
$ws.Range("L113").Value = 4890.125
$ws.Range("M113").Value = -299.9000000000001
$ws.Range("N113").Value = -9230.125
$ws.Range("H132").Value = 3563.3823
$ws.Range("I132").Value = 2141.0417
$ws.Range("K132").Value = 6423.125100000001
$ws.Range("M132").Value = -3893.125100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1149.5
$ws.Range("I82").Value = 1030
$ws.Range("J82").Value = 1388.5
$ws.Range("K82").Value = 1030
$ws.Range("L82").Value = 1388.5
$ws.Range("M82").Value = -669
$ws.Range("N82").Value = -2110.5
$ws.Range("H85").Value = 1149.5
$ws.Range("I85").Value = 1030
$ws.Range("J85").Value = 1388.5
$ws.Range("K85").Value = 1030
$ws.Range("L85").Value = 1388.5
$ws.Range("M85").Value = 218
$ws.Range("N85").Value = -3884.5
$ws.Range("H132").Value = 4920.3477
$ws.Range("I132").Value = 3722.2856
$ws.Range("K132").Value = 11166.8568
$ws.Range("M132").Value = -8636.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11631.25
$ws.Range("J45").Value = 11631.25
$ws.Range("L45").Value = 11631.25
$ws.Range("N45").Value = -12613.25
$ws.Range("H81").Value = 2230.0476
$ws.Range("I81").Value = 2040.7222
$ws.Range("K81").Value = 4081.4444
$ws.Range("M81").Value = -3020.4444
$ws.Range("H84").Value = 2230.0476
$ws.Range("I84").Value = 2040.7222
$ws.Range("K84").Value = 20407.222
$ws.Range("M84").Value = -15103.222
$ws.Range("H96").Value = 21489.846
$ws.Range("I96").Value = 2233
$ws.Range("J96").Value = 27266.9
$ws.Range("K96").Value = 2233
$ws.Range("L96").Value = 27266.9
$ws.Range("M96").Value = -860
$ws.Range("N96").Value = -30012.9
$ws.Range("H136").Value = 7641.6924
$ws.Range("I136").Value = 6914.931
$ws.Range("J136").Value = 9749.299999999999
$ws.Range("K136").Value = 20744.793
$ws.Range("L136").Value = 29247.9
$ws.Range("M136").Value = -18194.793
$ws.Range("N136").Value = -34347.89999999999
